$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.747.57"
$ws.Range("E2").Value = "  -2.01%  "

$ws.Range("D3").Value = "3.767.05"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'597.74"
$ws.Range("E5").Value = "  -2.77%  "

$ws.Range("D6").Value = "'175.92"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").Value = "3.764.67"
$ws.Range("E7").Value = "  +0.61%  "

$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  +0.39%  "

$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -3.76%  "

$ws.Range("E11").Value = "  -5.46%  "

$ws.Range("E12").Value = "  -3.78%  "

$ws.Range("D13").Value = "'38.62"
$ws.Range("E13").Value = "  -3.35%  "

$ws.Range("D14").Value = "'0.0000246"
$ws.Range("E14").Value = "  -2.80%  "

$ws.Range("D15").Value = "4.400.17"
$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").Value = "3.767.76"
$ws.Range("E16").Value = "  +0.68%  "

$ws.Range("D17").Value = "67.659.08"
$ws.Range("E17").Value = "  -2.28%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.21"
$ws.Range("E18").Value = "  -3.27%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.115"
$ws.Range("E19").Value = "  -4.22%  "

$ws.Range("D20").Value = "'16.52"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").Value = "'490.41"
$ws.Range("E21").Value = "  -1.57%  "

$ws.Range("E22").Value = "  -2.79%  "

$ws.Range("E23").Value = "  +2.77%  "

$ws.Range("E24").Value = "  +12.50%  "

$ws.Range("D25").Value = "'85.34"
$ws.Range("E25").Value = "  -0.54%  "

$ws.Range("E26").Value = "  -6.22%  "

$ws.Range("D27").Value = "'12.30"
$ws.Range("E27").Value = "  -3.93%  "

$ws.Range("D28").Value = "'10.23"
$ws.Range("E28").Value = "  -4.40%  "

$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("E30").Value = "  +0.72%  "

$ws.Range("E31").Value = "  -2.53%  "

$ws.Range("D32").Value = "'32.23"
$ws.Range("E32").Value = "  +5.64%  "

$ws.Range("D33").Value = "'7.76"
$ws.Range("E33").Value = "  -2.68%  "

$ws.Range("E34").Value = "  -3.64%  "

$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("E36").Value = "  -3.67%  "

$ws.Range("D37").Value = "'5.80"
$ws.Range("E37").Value = "  -4.77%  "

$ws.Range("E38").Value = "  -1.35%  "

$ws.Range("D39").Value = "'0.330"
$ws.Range("E39").Value = "  -4.89%  "

$ws.Range("D40").Value = "'448.71"
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("D41").Value = "'49.06"
$ws.Range("E41").Value = "  -1.31%  "

$ws.Range("D42").Value = "'2.01"
$ws.Range("E42").Value = "  -2.72%  "

$ws.Range("D43").Value = "'2.91"
$ws.Range("E43").Value = "  -3.47%  "

$ws.Range("D44").Value = "'8.37"
$ws.Range("E44").Value = "  -2.17%  "

$ws.Range("D45").Value = "'41.31"
$ws.Range("E45").Value = "  -7.62%  "

$ws.Range("D46").Value = "2.833.91"
$ws.Range("E46").Value = "  -3.58%  "

$ws.Range("D48").Value = "'138.67"
$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("D49").Value = "'0.0351"
$ws.Range("E49").Value = "  -2.22%  "

$ws.Range("D50").Value = "'26.18"
$ws.Range("E50").Value = "  -4.13%  "

$ws.Range("D51").Value = "'23.87"
$ws.Range("E51").Value = "  +8.35%  "
